$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 143 - this shifts the existing rows 143:209
# down to 144:210, matching the target dimension A1:T210.
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new market-price record.
$ws.Range("A143").Value = 5
$ws.Range("B143").Value = "Macroferia Regional de Talca"
$ws.Range("C143").Value = "Maule"
$ws.Range("D143").Value = 44510
$ws.Range("E143").Value = 7
$ws.Range("F143").Value = "Fruta"
$ws.Range("G143").Value = 100101
$ws.Range("H143").Value = "Berries"
$ws.Range("I143").Value = 100101007
$ws.Range("J143").Value = "Kiwi"
$ws.Range("K143").Value = "Hayward"
$ws.Range("L143").Value = "Primera"
$ws.Range("M143").Value = 160
$ws.Range("N143").Value = 18000
$ws.Range("O143").Value = 18000
$ws.Range("P143").Value = 18000
$ws.Range("Q143").Value = "$/bandeja 18 kilos"
$ws.Range("R143").Value = "Provincia de Curicó"
$ws.Range("S143").Value = 1000
$ws.Range("T143").Value = 18
